# NSMB - movie done!!
# Update final rows of the V4 (sheet1) splits table: fix B/C/D for rows 213-216
# and append new rows 217-227 for the remainder of the run (checkpoints,
# final door, end of input, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 213 ---
$ws.Range("A213").Value = "Checkpoint"
$ws.Range("B213").Value = 74749
$ws.Range("C213").Value = 86590

# --- Row 214 ---
$ws.Range("A214").Value = "Enter door"
$ws.Range("B214").Value = 74982
$ws.Range("C214").Value = 86830

# --- Row 215 ---
$ws.Range("A215").Value = "Speed = 2"
$ws.Range("B215").Value = 75107
$ws.Range("C215").Value = 86956

# --- Row 216 ---
$ws.Range("A216").Value = "Enter door"
$ws.Range("B216").Value = 75280
$ws.Range("C216").Value = 87129

# --- Row 217 (new) ---
$ws.Range("A217").Value = "Enter door"
$ws.Range("B217").Value = 75591
$ws.Range("C217").Value = 87443

# --- Row 218 (new) ---
$ws.Range("A218").Value = "Checkpoint 7/4"
$ws.Range("B218").Value = 75741
$ws.Range("C218").Value = 87604

# --- Row 219 (new) ---
$ws.Range("A219").Value = "Checkpoint 305/304"
$ws.Range("B219").Value = 75828
$ws.Range("C219").Value = 87692

# --- Row 220 (new) ---
$ws.Range("A220").Value = "Checkpoint 738/736"
$ws.Range("B220").Value = 75972
$ws.Range("C220").Value = 87836

# --- Row 221 (new) ---
$ws.Range("A221").Value = "Checkpoint 1505/1503"
$ws.Range("B221").Value = 76227
$ws.Range("C221").Value = 88091

# --- Row 222 (new) ---
$ws.Range("A222").Value = "Checkpoint 1742/1740"
$ws.Range("B222").Value = 76306
$ws.Range("C222").Value = 88170

# --- Row 223 (new) ---
$ws.Range("A223").Value = "Checkpoint 2158"
$ws.Range("B223").Value = 76444
$ws.Range("C223").Value = 88308

# --- Row 224 (new) ---
$ws.Range("A224").Value = "Enter Final Door"
$ws.Range("B224").Value = 76741
$ws.Range("C224").Value = 88605

# --- Row 225 (new) ---
$ws.Range("A225").Value = "Speed = 0"
$ws.Range("B225").Value = 77063
$ws.Range("C225").Value = 88946

# --- Row 226 (new) ---
$ws.Range("A226").Value = "End Input"
$ws.Range("B226").Value = 77862
$ws.Range("C226").Value = 89784

# --- Row 227 (new) ---
$ws.Range("A227").Value = "Touch Button"
$ws.Range("B227").Value = 77899
$ws.Range("C227").Value = 89784

# Fill the shared "Diff" formula down through the newly added rows so
# D90:D227 all compute IF(B>0,C-B,0) like the rest of the column.
for ($r = 213; $r -le 227; $r++) {
    $ws.Range("D$r").Formula = "=IF(B$r >  0,C$r-B$r, 0)"
}

# Leave the active cell on the row right after the last entry, matching
# where Excel drops you after typing the last row of data.
$ws.Range("B228").Select()

Write-Host "done"
